# Scheduled market-price refresh: update currentAveragePrice* / LevePrice* /
# LeveProfit* columns (H:N) for the affected leve rows across the Carbuncle
# profit-tracking sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 5
$ws.Range("H5").Value = 249.16667
$ws.Range("J5").Value = 114.666664
$ws.Range("L5").Value = 114.666664
$ws.Range("N5").Value = -344.666664

# Row 137
$ws.Range("H137").Value = 1126.5
$ws.Range("I137").Value = 1144.9524
$ws.Range("J137").Value = 1071.1428
$ws.Range("K137").Value = 3434.857199999999
$ws.Range("L137").Value = 3213.4284
$ws.Range("M137").Value = -884.8571999999995
$ws.Range("N137").Value = -8313.428400000001

# Row 138
$ws.Range("H138").Value = 3511.5
$ws.Range("I138").Value = 2394.6584
$ws.Range("J138").Value = 4287.6104
$ws.Range("K138").Value = 7183.975199999999
$ws.Range("L138").Value = 12862.8312
$ws.Range("M138").Value = -2043.975199999999
$ws.Range("N138").Value = -23142.8312

$ws = $wb.Worksheets.Item("ARM")
# Row 4
$ws.Range("H4").Value = 1000
$ws.Range("I4").Value = 1000
$ws.Range("K4").Value = 1000
$ws.Range("M4").Value = -884

# Row 9
$ws.Range("H9").Value = 50000
$ws.Range("J9").Value = 50000
$ws.Range("L9").Value = 50000
$ws.Range("N9").Value = -50340

# Row 20
$ws.Range("H20").Value = 50000
$ws.Range("J20").Value = 50000
$ws.Range("L20").Value = 50000
$ws.Range("N20").Value = -50540

# Row 32
$ws.Range("H32").Value = 10109.507
$ws.Range("I32").Value = 7501.9253
$ws.Range("J32").Value = 27580.3
$ws.Range("K32").Value = 7501.9253
$ws.Range("L32").Value = 27580.3
$ws.Range("M32").Value = -7214.9253
$ws.Range("N32").Value = -28154.3

# Row 37
$ws.Range("H37").Value = 10000
$ws.Range("I37").Value = 10000
$ws.Range("K37").Value = 10000
$ws.Range("M37").Value = -9727

# Row 44
$ws.Range("H44").Value = 28510
$ws.Range("J44").Value = 28510
$ws.Range("L44").Value = 28510
$ws.Range("N44").Value = -29486

# Row 55
$ws.Range("H55").Value = 30095
$ws.Range("J55").Value = 38793.332
$ws.Range("L55").Value = 38793.332
$ws.Range("N55").Value = -39423.332

# Row 61
$ws.Range("H61").Value = 1920.6552
$ws.Range("I61").Value = 1539.3889
$ws.Range("J61").Value = 2544.5454
$ws.Range("K61").Value = 1539.3889
$ws.Range("L61").Value = 2544.5454
$ws.Range("M61").Value = -1327.3889
$ws.Range("N61").Value = -2968.5454

# Row 74
$ws.Range("H74").Value = 1955.4517
$ws.Range("I74").Value = 1304.7727
$ws.Range("J74").Value = 3546
$ws.Range("K74").Value = 1304.7727
$ws.Range("L74").Value = 3546
$ws.Range("M74").Value = -430.7727
$ws.Range("N74").Value = -5294

# Row 77
$ws.Range("H77").Value = 1955.4517
$ws.Range("I77").Value = 1304.7727
$ws.Range("J77").Value = 3546
$ws.Range("K77").Value = 6523.863499999999
$ws.Range("L77").Value = 17730
$ws.Range("M77").Value = -2155.863499999999
$ws.Range("N77").Value = -26466

# Row 136
$ws.Range("H136").Value = 1920.6552
$ws.Range("I136").Value = 1539.3889
$ws.Range("J136").Value = 2544.5454
$ws.Range("K136").Value = 4618.1667
$ws.Range("L136").Value = 7633.6362
$ws.Range("M136").Value = -2068.1667
$ws.Range("N136").Value = -12733.6362

# Row 140
$ws.Range("H140").Value = 39190
$ws.Range("J140").Value = 39190
$ws.Range("L140").Value = 39190
$ws.Range("N140").Value = -49550

$ws = $wb.Worksheets.Item("BSM")
# Row 140
$ws.Range("H140").Value = 56655.715
$ws.Range("J140").Value = 56655.715
$ws.Range("L140").Value = 56655.715
$ws.Range("N140").Value = -67015.715

# Row 141
$ws.Range("H141").Value = 44938
$ws.Range("J141").Value = 44938
$ws.Range("L141").Value = 44938
$ws.Range("N141").Value = -55298

$ws = $wb.Worksheets.Item("CRP")
# Row 17
$ws.Range("H17").Value = 8000
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()

# Row 31
$ws.Range("H31").Value = 4251.101
$ws.Range("I31").Value = 2044.2106
$ws.Range("J31").Value = 8182.125
$ws.Range("K31").Value = 2044.2106
$ws.Range("L31").Value = 8182.125
$ws.Range("M31").Value = -1749.2106
$ws.Range("N31").Value = -8772.125

# Row 34
$ws.Range("H34").Value = 4251.101
$ws.Range("I34").Value = 2044.2106
$ws.Range("J34").Value = 8182.125
$ws.Range("K34").Value = 2044.2106
$ws.Range("L34").Value = 8182.125
$ws.Range("M34").Value = -1842.2106
$ws.Range("N34").Value = -8586.125

# Row 122
$ws.Range("H122").Value = 4964.25
$ws.Range("I122").Value = 1300
$ws.Range("J122").Value = 8628.5
$ws.Range("K122").Value = 3900
$ws.Range("L122").Value = 25885.5
$ws.Range("M122").Value = -1450
$ws.Range("N122").Value = -30785.5

# Row 134
$ws.Range("H134").Value = 2017.5625
$ws.Range("I134").Value = 1967.2858
$ws.Range("J134").Value = 2369.5
$ws.Range("K134").Value = 5901.857400000001
$ws.Range("L134").Value = 7108.5
$ws.Range("M134").Value = -3366.857400000001
$ws.Range("N134").Value = -12178.5

# Row 138
$ws.Range("H138").Value = 41204.332
$ws.Range("J138").Value = 41204.332
$ws.Range("L138").Value = 41204.332
$ws.Range("N138").Value = -51484.332

# Row 140
$ws.Range("H140").Value = 81646.164
$ws.Range("J140").Value = 81646.164
$ws.Range("L140").Value = 81646.164
$ws.Range("N140").Value = -92006.164

$ws = $wb.Worksheets.Item("CUL")
# Row 33
$ws.Range("H33").Value = 90.625
$ws.Range("J33").Value = 100
$ws.Range("N33").Value = -1166

# Row 107
$ws.Range("H107").Value = 358678.38
$ws.Range("J107").Value = 524059.03
$ws.Range("L107").Value = 1572177.09
$ws.Range("N107").Value = -1576017.09

# Row 113
$ws.Range("H113").Value = 671.2
$ws.Range("I113").Value = 714.2857
$ws.Range("J113").Value = 570.6667
$ws.Range("K113").Value = 2142.8571
$ws.Range("L113").Value = 1712.0001
$ws.Range("M113").Value = 27.14289999999983
$ws.Range("N113").Value = -6052.0001

# Row 122
$ws.Range("H122").Value = 1167.9
$ws.Range("I122").Value = 624.9375
$ws.Range("J122").Value = 3339.75
$ws.Range("K122").Value = 5624.4375
$ws.Range("L122").Value = 30057.75
$ws.Range("M122").Value = -3174.4375
$ws.Range("N122").Value = -34957.75

# Row 126
$ws.Range("H126").Value = 1461.5385
$ws.Range("I126").Value = 710
$ws.Range("J126").Value = 1598.1818
$ws.Range("K126").Value = 2130
$ws.Range("L126").Value = 4794.5454
$ws.Range("M126").Value = 2810
$ws.Range("N126").Value = -14674.5454

# Row 132
$ws.Range("H132").Value = 1918.1154
$ws.Range("I132").Value = 1088.2
$ws.Range("J132").Value = 2115.7144
$ws.Range("K132").Value = 9793.800000000001
$ws.Range("L132").Value = 19041.4296
$ws.Range("M132").Value = -7263.800000000001
$ws.Range("N132").Value = -24101.4296

# Row 137
$ws.Range("H137").Value = 6838.773
$ws.Range("I137").Value = 2241.9
$ws.Range("J137").Value = 8190.794
$ws.Range("K137").Value = 6725.700000000001
$ws.Range("L137").Value = 24572.382
$ws.Range("M137").Value = -1625.700000000001
$ws.Range("N137").Value = -34772.382

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 817.2
$ws.Range("I46").Value = 771.5
$ws.Range("J46").Value = 1000
$ws.Range("K46").Value = 771.5
$ws.Range("L46").Value = 1000
$ws.Range("M46").Value = -583.5
$ws.Range("N46").Value = -1376

# Row 58
$ws.Range("H58").Value = 6000
$ws.Range("I58").Value = 6000
$ws.Range("K58").Value = 6000
$ws.Range("M58").Value = -5740

# Row 68
$ws.Range("H68").Value = 1653.625
$ws.Range("I68").Value = 1198
$ws.Range("J68").Value = 1805.5
$ws.Range("K68").Value = 1198
$ws.Range("L68").Value = 1805.5
$ws.Range("M68").Value = -449
$ws.Range("N68").Value = -3303.5

# Row 71
$ws.Range("H71").Value = 1653.625
$ws.Range("I71").Value = 1198
$ws.Range("J71").Value = 1805.5
$ws.Range("K71").Value = 5990
$ws.Range("L71").Value = 9027.5
$ws.Range("M71").Value = -2246
$ws.Range("N71").Value = -16515.5

$ws = $wb.Worksheets.Item("WVR")
# Row 29
$ws.Range("H29").Value = 27450
$ws.Range("I29").Value = 4900
$ws.Range("J29").Value = 50000
$ws.Range("K29").Value = 4900
$ws.Range("L29").Value = 50000
$ws.Range("M29").Value = -4610
$ws.Range("N29").Value = -50580
